# Update "想去人数" (interested-people count) values in column F across sheets
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2821
$ws1.Range("F3").Value = 1140
$ws1.Range("F4").Value = 20618
$ws1.Range("F6").Value = 2633
$ws1.Range("F8").Value = 615
$ws1.Range("F10").Value = 737
$ws1.Range("F11").Value = 273
$ws1.Range("F12").Value = 258
$ws1.Range("F14").Value = 103
$ws1.Range("F17").Value = 245
$ws1.Range("F19").Value = 405
$ws1.Range("F20").Value = 15
$ws1.Range("F21").Value = 25

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 26
$ws2.Range("F5").Value = 319

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6090
$ws3.Range("F3").Value = 684
$ws3.Range("F4").Value = 658
$ws3.Range("F5").Value = 1461
$ws3.Range("F6").Value = 45

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6090
$ws4.Range("F3").Value = 684
$ws4.Range("F4").Value = 658
$ws4.Range("F5").Value = 1461
$ws4.Range("F6").Value = 2821
$ws4.Range("F7").Value = 1140
$ws4.Range("F8").Value = 20618
$ws4.Range("F10").Value = 26
$ws4.Range("F13").Value = 319
$ws4.Range("F14").Value = 2633
$ws4.Range("F17").Value = 45
$ws4.Range("F18").Value = 615
$ws4.Range("F20").Value = 737
$ws4.Range("F21").Value = 273
$ws4.Range("F22").Value = 258
$ws4.Range("F27").Value = 103
$ws4.Range("F34").Value = 245
$ws4.Range("F38").Value = 405
$ws4.Range("F40").Value = 15
$ws4.Range("F41").Value = 25
